# Auto-generated COM-interop script applying the "extraction_status" column split
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comentarios")

# Insert a new column at O (shifts old O:"created_time_raw" data into column P)
$ws.Columns.Item(15).Insert()

# New header for the inserted column, and fix style/number-format on column I (fecha_comentario)
$ws.Range("O1").Value = 'extraction_status'
$ws.Range("I2:I19").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Move the raw extractor JSON payload that used to live in column O into column P,
# and fix the likes_count cells (K10:K18) that were stored as text back to numbers.
$ws.Range("P2").Value = '{''videoWebUrl'': ''https://www.tiktok.com/@alpinacol/video/7574526148311289096'', ''submittedVideoUrl'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''input'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''cid'': ''7575394784510427922'', ''createTime'': 1763784064, ''createTimeISO'': ''2025-11-22T04:01:04.000Z'', ''text'': ''🥰 Que lindo que no es IA'', ''diggCount'': 0, ''likedByAuthor'': False, ''pinnedByAuthor'': False, ''repliesToId'': None, ''replyCommentTotal'': 0, ''uid'': ''7061039853035422725'', ''uniqu'
$ws.Range("P3").Value = '{''videoWebUrl'': ''https://www.tiktok.com/@alpinacol/video/7574526148311289096'', ''submittedVideoUrl'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''input'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''cid'': ''7575362101340324626'', ''createTime'': 1763776456, ''createTimeISO'': ''2025-11-22T01:54:16.000Z'', ''text'': ''por qué no explica qué clase de mascarilla de vino hacer cómo se prepara qué ingredientes le hecho'', ''diggCount'': 0, ''likedByAuthor'': False, ''pinnedByAuthor'': False, ''replie'
$ws.Range("P4").Value = '{''videoWebUrl'': ''https://www.tiktok.com/@alpinacol/video/7574526148311289096'', ''submittedVideoUrl'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''input'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''cid'': ''7575359261501408018'', ''createTime'': 1763775782, ''createTimeISO'': ''2025-11-22T01:43:02.000Z'', ''text'': ''ese es el mapa al one piece?'', ''diggCount'': 0, ''likedByAuthor'': False, ''pinnedByAuthor'': False, ''repliesToId'': None, ''replyCommentTotal'': 0, ''uid'': ''7178493082228442117'', ''u'
$ws.Range("P5").Value = '{''videoWebUrl'': ''https://www.tiktok.com/@alpinacol/video/7574526148311289096'', ''submittedVideoUrl'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''input'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''cid'': ''7575314301649797906'', ''createTime'': 1763765325, ''createTimeISO'': ''2025-11-21T22:48:45.000Z'', ''text'': ''Me gusta que no es IA'', ''diggCount'': 0, ''likedByAuthor'': False, ''pinnedByAuthor'': False, ''repliesToId'': None, ''replyCommentTotal'': 0, ''uid'': ''7476626177487537207'', ''uniqueId'
$ws.Range("P6").Value = '{''videoWebUrl'': ''https://www.tiktok.com/@alpinacol/video/7574526148311289096'', ''submittedVideoUrl'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''input'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''cid'': ''7575264519527105300'', ''createTime'': 1763753731, ''createTimeISO'': ''2025-11-21T19:35:31.000Z'', ''text'': ''Bendiciones 🥰'', ''diggCount'': 1, ''likedByAuthor'': False, ''pinnedByAuthor'': False, ''repliesToId'': None, ''replyCommentTotal'': 0, ''uid'': ''7392723576502043654'', ''uniqueId'': ''ana0'
$ws.Range("P7").Value = '{''videoWebUrl'': ''https://www.tiktok.com/@alpinacol/video/7574526148311289096'', ''submittedVideoUrl'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''input'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''cid'': ''7575251936485180178'', ''createTime'': 1763750972, ''createTimeISO'': ''2025-11-21T18:49:32.000Z'', ''text'': ''Que bello video gracias Alpina😊'', ''diggCount'': 1, ''likedByAuthor'': False, ''pinnedByAuthor'': False, ''repliesToId'': None, ''replyCommentTotal'': 0, ''uid'': ''7520408854468035602'','
$ws.Range("P8").Value = '{''videoWebUrl'': ''https://www.tiktok.com/@alpinacol/video/7574526148311289096'', ''submittedVideoUrl'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''input'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''cid'': ''7575188600700011271'', ''createTime'': 1763736191, ''createTimeISO'': ''2025-11-21T14:43:11.000Z'', ''text'': ''☺️☺️😊'', ''diggCount'': 0, ''likedByAuthor'': False, ''pinnedByAuthor'': False, ''repliesToId'': None, ''replyCommentTotal'': 0, ''uid'': ''7417587759336457222'', ''uniqueId'': ''user72252409'
$ws.Range("P9").Value = '{''videoWebUrl'': ''https://www.tiktok.com/@alpinacol/video/7574526148311289096'', ''submittedVideoUrl'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''input'': ''https://www.tiktok.com/@/video/7574526148311289096'', ''cid'': ''7574663024414688001'', ''createTime'': 1763613848, ''createTimeISO'': ''2025-11-20T04:44:08.000Z'', ''text'': ''6'', ''diggCount'': 2, ''likedByAuthor'': False, ''pinnedByAuthor'': False, ''repliesToId'': None, ''replyCommentTotal'': 0, ''uid'': ''7515116442964198418'', ''uniqueId'': ''excelmerysalazar'
$ws.Range("P10").Value = '{''facebookUrl'': ''https://www.facebook.com/100064867445065/posts/pfbid02qLKsWnTscnt4fxoxbp2H8DLmNMVurZcBpUBuwZEZzTwyX1YkYmYENbHa5UHu6Dg8l'', ''commentUrl'': ''https://www.facebook.com/alpina/posts/pfbid02eoWzQNzHoYKoC98YGYQmLBqznN2TMr5HX4ynYsYVThQemN3AJ7CW8VMq65vm1JcZl?comment_id=1631648168208831'', ''id'': ''Y29tbWVudDoxMjc5Nzc0NDQ0MTk0NzkyXzE2MzE2NDgxNjgyMDg4MzE='', ''feedbackId'': ''ZmVlZGJhY2s6MTI3OTc3NDQ0NDE5NDc5Ml8xNjMxNjQ4MTY4MjA4ODMx'', ''date'': ''2025-11-23T23:37:25.000Z'', ''text'': ''Pp0⁹9p9⁹'', ''profileU'
$ws.Range("P11").Value = '{''facebookUrl'': ''https://www.facebook.com/100064867445065/posts/pfbid02qLKsWnTscnt4fxoxbp2H8DLmNMVurZcBpUBuwZEZzTwyX1YkYmYENbHa5UHu6Dg8l'', ''commentUrl'': ''https://www.facebook.com/alpina/posts/pfbid02eoWzQNzHoYKoC98YGYQmLBqznN2TMr5HX4ynYsYVThQemN3AJ7CW8VMq65vm1JcZl?comment_id=848516504333779'', ''id'': ''Y29tbWVudDoxMjc5Nzc0NDQ0MTk0NzkyXzg0ODUxNjUwNDMzMzc3OQ=='', ''feedbackId'': ''ZmVlZGJhY2s6MTI3OTc3NDQ0NDE5NDc5Ml84NDg1MTY1MDQzMzM3Nzk='', ''date'': ''2025-11-21T23:26:50.000Z'', ''text'': ''Hola te ves bn tienes'
$ws.Range("P12").Value = '{''facebookUrl'': ''https://www.facebook.com/100064867445065/posts/pfbid02qLKsWnTscnt4fxoxbp2H8DLmNMVurZcBpUBuwZEZzTwyX1YkYmYENbHa5UHu6Dg8l'', ''commentUrl'': ''https://www.facebook.com/alpina/posts/pfbid02eoWzQNzHoYKoC98YGYQmLBqznN2TMr5HX4ynYsYVThQemN3AJ7CW8VMq65vm1JcZl?comment_id=789278630809884'', ''id'': ''Y29tbWVudDoxMjc5Nzc0NDQ0MTk0NzkyXzc4OTI3ODYzMDgwOTg4NA=='', ''feedbackId'': ''ZmVlZGJhY2s6MTI3OTc3NDQ0NDE5NDc5Ml83ODkyNzg2MzA4MDk4ODQ='', ''date'': ''2025-11-21T06:15:56.000Z'', ''text'': ''Tinga linga linga tin'
$ws.Range("P13").Value = '{''facebookUrl'': ''https://www.facebook.com/100064867445065/posts/pfbid02qLKsWnTscnt4fxoxbp2H8DLmNMVurZcBpUBuwZEZzTwyX1YkYmYENbHa5UHu6Dg8l'', ''commentUrl'': ''https://www.facebook.com/alpina/posts/pfbid02eoWzQNzHoYKoC98YGYQmLBqznN2TMr5HX4ynYsYVThQemN3AJ7CW8VMq65vm1JcZl?comment_id=703351456156856'', ''id'': ''Y29tbWVudDoxMjc5Nzc0NDQ0MTk0NzkyXzcwMzM1MTQ1NjE1Njg1Ng=='', ''feedbackId'': ''ZmVlZGJhY2s6MTI3OTc3NDQ0NDE5NDc5Ml83MDMzNTE0NTYxNTY4NTY='', ''date'': ''2025-11-19T21:27:18.000Z'', ''text'': ''Yamen l al policía pa'
$ws.Range("P14").Value = '{''facebookUrl'': ''https://www.facebook.com/100064867445065/posts/pfbid036iaCTpHgEiiR2DT4ZQ8FoZtCogRjyg3nAXxtVXh3XXRdduTJDg4fThnwZ742KgvPl'', ''commentUrl'': ''https://www.facebook.com/reel/2001022834075363/?comment_id=1011405657841959'', ''id'': ''Y29tbWVudDoxMjc5NzgxMDMwODYwODAwXzEwMTE0MDU2NTc4NDE5NTk='', ''feedbackId'': ''ZmVlZGJhY2s6MTI3OTc4MTAzMDg2MDgwMF8xMDExNDA1NjU3ODQxOTU5'', ''date'': ''2025-11-20T12:45:06.000Z'', ''text'': ''Listo'', ''profileUrl'': ''https://www.facebook.com/esperanza.londono.35325'', ''profileP'
$ws.Range("P15").Value = '{''facebookUrl'': ''https://www.facebook.com/100064867445065/posts/pfbid036iaCTpHgEiiR2DT4ZQ8FoZtCogRjyg3nAXxtVXh3XXRdduTJDg4fThnwZ742KgvPl'', ''commentUrl'': ''https://www.facebook.com/reel/2001022834075363/?comment_id=1201407255170646'', ''id'': ''Y29tbWVudDoxMjc5NzgxMDMwODYwODAwXzEyMDE0MDcyNTUxNzA2NDY='', ''feedbackId'': ''ZmVlZGJhY2s6MTI3OTc4MTAzMDg2MDgwMF8xMjAxNDA3MjU1MTcwNjQ2'', ''date'': ''2025-11-20T11:13:11.000Z'', ''text'': ''ZAMBRANO UNA  VERGUENSA!!! APROVECHADO Y AMBICIOSO!! Q PENA!!!!'', ''profileUrl'': ''ht'
$ws.Range("P16").Value = '{''facebookUrl'': ''https://www.facebook.com/100064672685926/posts/1258767229622355/'', ''commentUrl'': ''https://www.facebook.com/AlpinaBonYurt/posts/pfbid02vDjJopkwJ1E4ttEJPs8Q5WtCbkyYAnJ22tMuNREZgXTVE7PAeoMpaduR3LMCgVnAl?comment_id=861351922903088'', ''id'': ''Y29tbWVudDoxMjU4NzY3NzY5NjIyMzAxXzg2MTM1MTkyMjkwMzA4OA=='', ''feedbackId'': ''ZmVlZGJhY2s6MTI1ODc2Nzc2OTYyMjMwMV84NjEzNTE5MjI5MDMwODg='', ''date'': ''2025-11-25T05:11:04.000Z'', ''text'': ''Esq eso era muy caro un riñón prácticamente ps pero vale la pena'', ''p'
$ws.Range("P17").Value = '{''facebookUrl'': ''https://www.facebook.com/100064672685926/posts/1258767229622355/'', ''commentUrl'': ''https://www.facebook.com/AlpinaBonYurt/posts/pfbid02vDjJopkwJ1E4ttEJPs8Q5WtCbkyYAnJ22tMuNREZgXTVE7PAeoMpaduR3LMCgVnAl?comment_id=1820814188796053'', ''id'': ''Y29tbWVudDoxMjU4NzY3NzY5NjIyMzAxXzE4MjA4MTQxODg3OTYwNTM='', ''feedbackId'': ''ZmVlZGJhY2s6MTI1ODc2Nzc2OTYyMjMwMV8xODIwODE0MTg4Nzk2MDUz'', ''date'': ''2025-11-25T04:34:17.000Z'', ''text'': ''Cuando eran ricos y traían más que aire'', ''profileUrl'': ''https://www'
$ws.Range("P18").Value = '{''facebookUrl'': ''https://www.facebook.com/100064672685926/posts/1258767229622355/'', ''commentUrl'': ''https://www.facebook.com/AlpinaBonYurt/posts/pfbid02vDjJopkwJ1E4ttEJPs8Q5WtCbkyYAnJ22tMuNREZgXTVE7PAeoMpaduR3LMCgVnAl?comment_id=843566781589283'', ''id'': ''Y29tbWVudDoxMjU4NzY3NzY5NjIyMzAxXzg0MzU2Njc4MTU4OTI4Mw=='', ''feedbackId'': ''ZmVlZGJhY2s6MTI1ODc2Nzc2OTYyMjMwMV84NDM1NjY3ODE1ODkyODM='', ''date'': ''2025-11-25T03:23:08.000Z'', ''text'': ''Pos si como tres me da un coma diabético....'', ''profileUrl'': ''https:/'
$ws.Range("P19").Value = '{''postUrl'': ''https://www.instagram.com/p/DRPexY8jPUB/'', ''commentUrl'': ''https://www.instagram.com/p/DRPexY8jPUB/c/18364850800091592'', ''id'': ''18364850800091592'', ''text'': ''De la Moradita'', ''ownerUsername'': ''angelodominic61'', ''ownerProfilePicUrl'': ''https://scontent-sea5-1.cdninstagram.com/v/t51.2885-19/451569933_490434486865115_7000517521910309172_n.jpg?stp=dst-jpg_e0_s150x150_tt6&_nc_cat=109&ig_cache_key=GA1p6hrbDGAtDL4BADShXImY1SZhbkULAAAB1501500j-ccb7-5&ccb=7-5&_nc_sid=669407&efg=eyJ2ZW5jb2RlX3Rh'

# K10:K18 were text "0"/"1" in the source export; restore numeric types to match the rest of the column
$ws.Range("K10").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("K12").Value = 1
$ws.Range("K13").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("K18").Value = 0
